# Logbog.docx - "Sidste update af Logbog 04-10-2019"
#
# Append two new bullet-style log lines after the "Sprint 2 (4) går
# Namkinh i gang med." paragraph, matching the indentation/size of the
# surrounding entries, and move the hidden "_GoBack" bookmark so that it
# ends up right after the very last of the newly added text (mirroring
# where Word leaves it after the last edit).

$d = $word.ActiveDocument

$targetText = "Sprint 2 (4) går Namkinh i gang med."

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text.TrimEnd() -eq $targetText) {
        $target = $candidate
        break
    }
}

$firstNewText  = "Danskificeret Slime, så den har et dansk flag skin på sig."
$secondNewText = "Gjort mentalt klar til tilføjelsen af scoreboard"

# Insert a new (empty) paragraph right after the target paragraph; it
# naturally inherits the paragraph/run formatting (ind left=360, sz=36).
$r = $target.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

$p1 = $target.Next()
$p1.Range.Text = $firstNewText

$r1 = $p1.Range
$r1.Collapse(0)
$r1.InsertParagraphAfter()

$p2 = $p1.Next()

# Write the second paragraph's text with one extra sentinel character so
# that the position right after the real text is never the *last*
# character offset inside the paragraph while we create the zero-length
# bookmark range (a transient, non-boundary offset is used, then the
# sentinel is removed afterwards).
$p2.Range.Text = $secondNewText + "#"

$textEnd = $p2.Range.Start + $secondNewText.Length

$old = $d.Bookmarks.Item("_GoBack")
$old.Delete()

$bmRange = $d.Range($textEnd, $textEnd)
$d.Bookmarks.Add("_GoBack", $bmRange)

$sentinel = $d.Range($textEnd, $textEnd + 1)
$sentinel.Delete()
